$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.127.99"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "'1.656.17"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'217.94"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").Value = "'0.5262"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'0.2607"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").Value = "'0.06350"
$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("D10").Value = "'20.47"
$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").Value = "'4.511"
$ws.Range("E12").Value = "  +2.22%  "

$ws.Range("D13").Value = "'1.632.48"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").Value = "'0.5486"
$ws.Range("E14").Value = "  +0.96%  "

$ws.Range("D15").Value = "'0.0₅8239"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "'65.40"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "'26.140.93"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D19").Value = "'4.576"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").Value = "'192.15"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "'6.042"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").Value = "'141.82"
$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").Value = "'0.1251"
$ws.Range("E25").Value = "  +2.40%  "

$ws.Range("D26").Value = "'7.275"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("D29").Value = "'0.05922"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("D30").Value = "'1.278"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "'3.257"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").Value = "'1.586"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").Value = "'0.9534"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "'0.5702"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").Value = "'0.01620"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").Value = "'0.8484"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").Value = "'5.780"
$ws.Range("E40").Value = "  -3.61%  "

$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.030.42"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'103.02"
$ws.Range("E43").Value = "  +2.88%  "

$ws.Range("D44").Value = "'1.802.18"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'57.44"
$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").Value = "'0.4302"
$ws.Range("E47").Value = "  +2.07%  "

$ws.Range("D48").Value = "'1.482"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").Value = "'0.05164"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("D50").Value = "'7.823"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").Value = "'0.09727"
$ws.Range("E51").Value = "  +0.47%  "
